$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column C (Lawyers) for the rows we are about to fill as Text,
# so the numeric-looking "1"/"0" values are stored as text (matching the
# existing "str" cell type) instead of being coerced to numbers. The
# "General" format is restored below so the cell style index is unaffected.
$ws.Range("C2:C38").NumberFormat = "@"

$ws.Range("A2").Value = 'Huiye Law'
$ws.Range("B2").Value = '32sec'
$ws.Range("C2").Value = '1'
$ws.Range("A3").Value = 'Kinstellar'
$ws.Range("B3").Value = '34sec'
$ws.Range("C3").Value = '1'
$ws.Range("A4").Value = 'Njord Law'
$ws.Range("B4").Value = '48sec'
$ws.Range("C4").Value = '1'
$ws.Range("A5").Value = 'Keystone Law'
$ws.Range("B5").Value = '2min 29sec'
$ws.Range("C5").Value = '1'
$ws.Range("A6").Value = 'Kromann Reumert'
$ws.Range("B6").Value = '16sec'
$ws.Range("C6").Value = '1'
$ws.Range("A7").Value = 'Pedersoli'
$ws.Range("B7").Value = '9sec'
$ws.Range("C7").Value = '1'
$ws.Range("A8").Value = 'Magnusson Law'
$ws.Range("B8").Value = '13sec'
$ws.Range("C8").Value = '1'
$ws.Range("A9").Value = 'Hannes Snellman'
$ws.Range("B9").Value = '21sec'
$ws.Range("C9").Value = '1'
$ws.Range("A10").Value = 'Arthur Cox'
$ws.Range("B10").Value = '21sec'
$ws.Range("C10").Value = '1'
$ws.Range("A11").Value = 'Spencer West'
$ws.Range("B11").Value = '21sec'
$ws.Range("C11").Value = '1'
$ws.Range("A12").Value = 'Ellex'
$ws.Range("B12").Value = '24sec'
$ws.Range("C12").Value = '1'
$ws.Range("A13").Value = 'Deacons'
$ws.Range("B13").Value = '13sec'
$ws.Range("C13").Value = '1'
$ws.Range("A14").Value = 'Clemens Law'
$ws.Range("B14").Value = '7sec'
$ws.Range("C14").Value = '1'
$ws.Range("A15").Value = 'Remfry And Sagar'
$ws.Range("B15").Value = '44sec'
$ws.Range("C15").Value = '1'
$ws.Range("A16").Value = 'Walkers'
$ws.Range("B16").Value = '44sec'
$ws.Range("C16").Value = '1'
$ws.Range("A17").Value = 'DGKV'
$ws.Range("B17").Value = '25sec'
$ws.Range("C17").Value = '1'
$ws.Range("A18").Value = 'Conyers'
$ws.Range("B18").Value = '23sec'
$ws.Range("C18").Value = '1'
$ws.Range("A19").Value = 'Fox And Mandal'
$ws.Range("B19").Value = '15sec'
$ws.Range("C19").Value = '0'
$ws.Range("A20").Value = 'Ecija'
$ws.Range("B20").Value = '34sec'
$ws.Range("C20").Value = '1'
$ws.Range("A21").Value = 'Meitar Law Offices'
$ws.Range("B21").Value = '19sec'
$ws.Range("C21").Value = '1'
$ws.Range("A22").Value = 'Szecskay'
$ws.Range("B22").Value = '17sec'
$ws.Range("C22").Value = '1'
$ws.Range("A23").Value = 'KRB Law Firm'
$ws.Range("B23").Value = '15sec'
$ws.Range("C23").Value = '1'
$ws.Range("A24").Value = 'Brigrard Urrutia'
$ws.Range("B24").Value = '12sec'
$ws.Range("C24").Value = '1'
$ws.Range("A25").Value = 'Myers Fletcher And Gordon'
$ws.Range("B25").Value = '15sec'
$ws.Range("C25").Value = '1'
$ws.Range("A26").Value = 'Fischer'
$ws.Range("B26").Value = '18sec'
$ws.Range("C26").Value = '1'
$ws.Range("A27").Value = 'Spruson And Ferguson'
$ws.Range("B27").Value = '10sec'
$ws.Range("C27").Value = '1'
$ws.Range("A28").Value = 'Consortium Legal'
$ws.Range("B28").Value = '11sec'
$ws.Range("C28").Value = '1'
$ws.Range("A29").Value = 'Collas Crill'
$ws.Range("B29").Value = '9sec'
$ws.Range("C29").Value = '1'
$ws.Range("A30").Value = 'AL Goodbody'
$ws.Range("B30").Value = '4min 24sec'
$ws.Range("C30").Value = '1'
$ws.Range("A31").Value = 'Gitti And Partners Law Firm'
$ws.Range("B31").Value = '15sec'
$ws.Range("C31").Value = '1'
$ws.Range("A32").Value = 'Zepos And Yannopoulos'
$ws.Range("B32").Value = '14sec'
$ws.Range("C32").Value = '1'
$ws.Range("A33").Value = 'Carey Olsen'
$ws.Range("B33").Value = '55sec'
$ws.Range("C33").Value = '1'
$ws.Range("A34").Value = 'DahlLaw'
$ws.Range("B34").Value = '15sec'
$ws.Range("C34").Value = '1'
$ws.Range("A35").Value = 'Campbells Legal'
$ws.Range("B35").Value = '4sec'
$ws.Range("C35").Value = '1'
$ws.Range("A36").Value = 'Harneys'
$ws.Range("B36").Value = '2min 13sec'
$ws.Range("C36").Value = '1'
$ws.Range("A37").Value = 'Borenius'
$ws.Range("B37").Value = '6min 57sec'
$ws.Range("C37").Value = '1'
$ws.Range("A38").Value = 'Dechert LLP'
$ws.Range("B38").Value = '5sec'
$ws.Range("C38").Value = '1'

# Restore the default "General" number format on column C so the cell style
# index matches the rest of the sheet (only the stored value/type changes).
$ws.Range("C2:C38").NumberFormat = "General"

